$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split column D ("varejo.servico") into two columns:
#   D -> "servicos", E (new) -> "varejo"
$ws.Columns("E:E").Insert()

$ws.Range("D1").Value = "servicos"
$ws.Range("E1").Value = "varejo"

# Match column E's width to its neighbours (~13.57 chars, same as D/F)
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# New empty comment cell for "varejo" gets an underlined, wrap-text style
# (placeholder formatting, to be filled in later)
$e2 = $ws.Range("E2")
$e2.Font.Underline = 2
$e2.WrapText = $true

# Update the view: no more frozen/scrolled topLeftCell, selection now on E2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E2").Select()
